# Applies the "New crime data collected" update to the 111th Precinct
# weekly CompStat workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text updates (rich-text cells: only part of the text changes)
# ---------------------------------------------------------------------

# A8: "Volume 32   Number  38" -> "...Number  40"
$cell = $ws.Cells.Item(8, 1)
$txt = $cell.Value()
$old = "38"
$new = "40"
$pos = $txt.LastIndexOf($old)
$cell.Characters($pos + 1, $old.Length).Text = $new

# C9: "Report Covering the Week  9/15/2025  Through  9/21/2025"
#  -> "Report Covering the Week  9/29/2025  Through  10/5/2025"
$cell = $ws.Cells.Item(9, 3)
$txt = $cell.Value()
$old1 = "9/15/2025"
$new1 = "9/29/2025"
$pos1 = $txt.IndexOf($old1)
$cell.Characters($pos1 + 1, $old1.Length).Text = $new1

$txt = $cell.Value()
$old2 = "9/21/2025"
$new2 = "10/5/2025"
$pos2 = $txt.IndexOf($old2)
$cell.Characters($pos2 + 1, $old2.Length).Text = $new2

# ---------------------------------------------------------------------
# 2. Column E got a bit wider (auto best-fit after new data came in)
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 6.7

# ---------------------------------------------------------------------
# 3. Cells that flip between a numeric style and the "no data" text
#    style (shared strings "0" / "***.*") need their style copied from
#    a stable reference cell so the workbook reuses the existing style
#    entries instead of inventing new ones.
# ---------------------------------------------------------------------

# Reference cells (row 14 is untouched by this update):
#   C14 -> style "0"   (text, numFmt General)
#   E14 -> style "***.*" (text, numFmt General)
#   I14 -> plain integer style
#   K14 -> plain decimal style

# D16: 4 -> "0"
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(16, 4))
# E16: -25 -> "***.*"
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(16, 5))
# C17: 2 -> "0"
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(17, 3))
# F33: 1 -> "0"
$ws.Cells.Item(14, 3).Copy($ws.Cells.Item(33, 6))

# D28: "0" -> 1 (numeric again)
$ws.Cells.Item(14, 9).Copy($ws.Cells.Item(28, 4))
$ws.Cells.Item(28, 4).Value = 1
# E28: "***.*" -> 0 (numeric again)
$ws.Cells.Item(14, 11).Copy($ws.Cells.Item(28, 5))
$ws.Cells.Item(28, 5).Value = 0

# ---------------------------------------------------------------------
# 4. Plain numeric value updates (style unchanged)
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 6).Value = 3  # F15
$ws.Cells.Item(15, 8).Value = 200  # H15
$ws.Cells.Item(15, 9).Value = 10  # I15
$ws.Cells.Item(15, 10).Value = 8  # J15
$ws.Cells.Item(15, 11).Value = 25  # K15
$ws.Cells.Item(15, 12).Value = 11.111111111111  # L15
$ws.Cells.Item(16, 3).Value = 1  # C16
$ws.Cells.Item(16, 7).Value = 6  # G16
$ws.Cells.Item(16, 8).Value = 16.666666666666  # H16
$ws.Cells.Item(16, 9).Value = 33  # I16
$ws.Cells.Item(16, 10).Value = 48  # J16
$ws.Cells.Item(16, 11).Value = -31.25  # K16
$ws.Cells.Item(16, 12).Value = -35.294117647058  # L16
$ws.Cells.Item(16, 13).Value = -46.774193548387  # M16
$ws.Cells.Item(16, 14).Value = -86.134453781512  # N16
$ws.Cells.Item(17, 5).Value = -100  # E17
$ws.Cells.Item(17, 7).Value = 5  # G17
$ws.Cells.Item(17, 8).Value = 80  # H17
$ws.Cells.Item(17, 9).Value = 98  # I17
$ws.Cells.Item(17, 10).Value = 69  # J17
$ws.Cells.Item(17, 11).Value = 42.028985507246  # K17
$ws.Cells.Item(17, 12).Value = 34.246575342465  # L17
$ws.Cells.Item(17, 13).Value = 151.282051282051  # M17
$ws.Cells.Item(17, 14).Value = 7.692307692307  # N17
$ws.Cells.Item(18, 3).Value = 8  # C18
$ws.Cells.Item(18, 4).Value = 9  # D18
$ws.Cells.Item(18, 5).Value = -11.111111111111  # E18
$ws.Cells.Item(18, 6).Value = 29  # F18
$ws.Cells.Item(18, 7).Value = 28  # G18
$ws.Cells.Item(18, 8).Value = 3.571428571428  # H18
$ws.Cells.Item(18, 9).Value = 229  # I18
$ws.Cells.Item(18, 10).Value = 227  # J18
$ws.Cells.Item(18, 11).Value = 0.881057268722  # K18
$ws.Cells.Item(18, 12).Value = -1.716738197424  # L18
$ws.Cells.Item(18, 13).Value = 24.456521739130  # M18
$ws.Cells.Item(18, 14).Value = -71.085858585858  # N18
$ws.Cells.Item(19, 3).Value = 14  # C19
$ws.Cells.Item(19, 4).Value = 6  # D19
$ws.Cells.Item(19, 5).Value = 133.333333333333  # E19
$ws.Cells.Item(19, 6).Value = 56  # F19
$ws.Cells.Item(19, 7).Value = 29  # G19
$ws.Cells.Item(19, 8).Value = 93.103448275862  # H19
$ws.Cells.Item(19, 9).Value = 383  # I19
$ws.Cells.Item(19, 10).Value = 356  # J19
$ws.Cells.Item(19, 11).Value = 7.584269662921  # K19
$ws.Cells.Item(19, 12).Value = -20.867768595041  # L19
$ws.Cells.Item(19, 13).Value = 39.272727272727  # M19
$ws.Cells.Item(19, 14).Value = -8.373205741626  # N19
$ws.Cells.Item(20, 3).Value = 8  # C20
$ws.Cells.Item(20, 4).Value = 8  # D20
$ws.Cells.Item(20, 5).Value = 0  # E20
$ws.Cells.Item(20, 6).Value = 36  # F20
$ws.Cells.Item(20, 8).Value = 56.521739130434  # H20
$ws.Cells.Item(20, 9).Value = 244  # I20
$ws.Cells.Item(20, 10).Value = 199  # J20
$ws.Cells.Item(20, 11).Value = 22.613065326633  # K20
$ws.Cells.Item(20, 12).Value = 73.049645390070  # L20
$ws.Cells.Item(20, 13).Value = 121.818181818182  # M20
$ws.Cells.Item(20, 14).Value = -90.105433901054  # N20
$ws.Cells.Item(21, 3).Value = 32  # C21
$ws.Cells.Item(21, 4).Value = 25  # D21
$ws.Cells.Item(21, 5).Value = 28  # E21
$ws.Cells.Item(21, 6).Value = 140  # F21
$ws.Cells.Item(21, 7).Value = 92  # G21
$ws.Cells.Item(21, 8).Value = 52.173913043478  # H21
$ws.Cells.Item(21, 9).Value = 999  # I21
$ws.Cells.Item(21, 10).Value = 908  # J21
$ws.Cells.Item(21, 11).Value = 10.022026431718  # K21
$ws.Cells.Item(21, 12).Value = 0.604229607250  # L21
$ws.Cells.Item(21, 13).Value = 48  # M21
$ws.Cells.Item(21, 14).Value = -75.136884021901  # N21
$ws.Cells.Item(24, 3).Value = 7  # C24
$ws.Cells.Item(24, 4).Value = 10  # D24
$ws.Cells.Item(24, 5).Value = -30  # E24
$ws.Cells.Item(24, 6).Value = 41  # F24
$ws.Cells.Item(24, 7).Value = 49  # G24
$ws.Cells.Item(24, 8).Value = -16.326530612244  # H24
$ws.Cells.Item(24, 9).Value = 423  # I24
$ws.Cells.Item(24, 10).Value = 422  # J24
$ws.Cells.Item(24, 11).Value = 0.236966824644  # K24
$ws.Cells.Item(24, 12).Value = -4.081632653061  # L24
$ws.Cells.Item(24, 13).Value = 25.892857142857  # M24
$ws.Cells.Item(25, 4).Value = 3  # D25
$ws.Cells.Item(25, 5).Value = -33.333333333333  # E25
$ws.Cells.Item(25, 6).Value = 8  # F25
$ws.Cells.Item(25, 7).Value = 12  # G25
$ws.Cells.Item(25, 8).Value = -33.333333333333  # H25
$ws.Cells.Item(25, 9).Value = 78  # I25
$ws.Cells.Item(25, 10).Value = 89  # J25
$ws.Cells.Item(25, 11).Value = -12.359550561797  # K25
$ws.Cells.Item(25, 12).Value = -16.129032258064  # L25
$ws.Cells.Item(26, 3).Value = 7  # C26
$ws.Cells.Item(26, 4).Value = 3  # D26
$ws.Cells.Item(26, 5).Value = 133.333333333333  # E26
$ws.Cells.Item(26, 6).Value = 20  # F26
$ws.Cells.Item(26, 7).Value = 17  # G26
$ws.Cells.Item(26, 8).Value = 17.647058823529  # H26
$ws.Cells.Item(26, 9).Value = 166  # I26
$ws.Cells.Item(26, 10).Value = 166  # J26
$ws.Cells.Item(26, 11).Value = 0  # K26
$ws.Cells.Item(26, 12).Value = -1.775147928994  # L26
$ws.Cells.Item(26, 13).Value = 29.6875  # M26
$ws.Cells.Item(27, 6).Value = 3  # F27
$ws.Cells.Item(27, 8).Value = 200  # H27
$ws.Cells.Item(27, 9).Value = 13  # I27
$ws.Cells.Item(27, 10).Value = 11  # J27
$ws.Cells.Item(27, 11).Value = 18.181818181818  # K27
$ws.Cells.Item(27, 12).Value = 30  # L27
$ws.Cells.Item(28, 6).Value = 3  # F28
$ws.Cells.Item(28, 7).Value = 3  # G28
$ws.Cells.Item(28, 8).Value = 0  # H28
$ws.Cells.Item(28, 9).Value = 11  # I28
$ws.Cells.Item(28, 10).Value = 9  # J28
$ws.Cells.Item(28, 11).Value = 22.222222222222  # K28
$ws.Cells.Item(28, 12).Value = -21.428571428571  # L28
